$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 86.29678392075563
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 112.9448213429104
